# Adds two new columns, I ("I0") and J ("IF"), to Sheet1, mirroring the
# existing header/style formatting used by the other columns (A:H) and
# filling in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Match the style already applied to the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("I1").Style = $ws.Range("H1").Style

$ws.Range("J1").Value = "IF"
$ws.Range("J1").Style = $ws.Range("H1").Style

# --- Data rows (rows 2-39) ---------------------------------------------
$values = @(
    @(2, 6, 7),
    @(3, 6, 6),
    @(4, 9, 9),
    @(5, 7, 8),
    @(6, 8, 8),
    @(7, 9, 9),
    @(8, 7, 8),
    @(9, 8, 8),
    @(10, 10, 10),
    @(11, 9, 9),
    @(12, 7, 8),
    @(13, 3, 4),
    @(14, 6, 8),
    @(15, 8, 8),
    @(16, 8, 8),
    @(17, 8, 9),
    @(18, 7, 8),
    @(19, 4, 5),
    @(20, 6, 6),
    @(21, 6, 8),
    @(22, 9, 9),
    @(23, 10, 10),
    @(24, 8, 8),
    @(25, 9, 9),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 7, 7),
    @(29, 8, 9),
    @(30, 4, 5),
    @(31, 6, 6),
    @(32, 6, 6),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 5, 5),
    @(36, 5, 6),
    @(37, 5, 5),
    @(38, 4, 4),
    @(39, 6, 6)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]

    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
